$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(11, 8).Value = 99.8
$ws.Cells.Item(11, 9).Value = 99.8
$ws.Cells.Item(11, 11).Value = 99.8
$ws.Cells.Item(11, 13).Value = 40.2
$ws.Cells.Item(18, 8).Value = 3119.0833
$ws.Cells.Item(18, 9).Value = 3354.9075
$ws.Cells.Item(18, 10).Value = 996.6667
$ws.Cells.Item(18, 11).Value = 3354.9075
$ws.Cells.Item(18, 12).Value = 996.6667
$ws.Cells.Item(18, 13).Value = -3070.9075
$ws.Cells.Item(18, 14).Value = -1564.6667
$ws.Cells.Item(32, 8).Value = 2183.9412
$ws.Cells.Item(32, 10).Value = 1271.909
$ws.Cells.Item(32, 12).Value = 1271.909
$ws.Cells.Item(32, 14).Value = -1923.909
$ws.Cells.Item(43, 8).Value = 4154.3335
$ws.Cells.Item(43, 10).Value = 3655.7144
$ws.Cells.Item(43, 12).Value = 3655.7144
$ws.Cells.Item(43, 14).Value = -3793.7144
$ws.Cells.Item(62, 8).Value = 22925.572
$ws.Cells.Item(62, 9).Value = 16947
$ws.Cells.Item(62, 11).Value = 16947
$ws.Cells.Item(62, 13).Value = -16323
$ws.Cells.Item(65, 8).Value = 22925.572
$ws.Cells.Item(65, 9).Value = 16947
$ws.Cells.Item(65, 11).Value = 84735
$ws.Cells.Item(65, 13).Value = -81615
$ws.Cells.Item(86, 8).Value = 38366.918
$ws.Cells.Item(86, 9).Value = 1180.6
$ws.Cells.Item(86, 11).Value = 1180.6
$ws.Cells.Item(86, 13).Value = -57.59999999999991
$ws.Cells.Item(89, 8).Value = 38366.918
$ws.Cells.Item(89, 9).Value = 1180.6
$ws.Cells.Item(89, 11).Value = 5903
$ws.Cells.Item(89, 13).Value = -287
$ws.Cells.Item(98, 8).Value = 46309.215
$ws.Cells.Item(98, 9).Value = 61816
$ws.Cells.Item(98, 10).Value = 18397
$ws.Cells.Item(98, 11).Value = 61816
$ws.Cells.Item(98, 12).Value = 18397
$ws.Cells.Item(98, 13).Value = -60318
$ws.Cells.Item(98, 14).Value = -21393
$ws.Cells.Item(112, 8).Value = 56572.15
$ws.Cells.Item(112, 10).Value = 86337.92
$ws.Cells.Item(112, 12).Value = 259013.76
$ws.Cells.Item(112, 14).Value = -261229.76
$ws.Cells.Item(122, 8).Value = 46309.215
$ws.Cells.Item(122, 9).Value = 61816
$ws.Cells.Item(122, 10).Value = 18397
$ws.Cells.Item(122, 11).Value = 185448
$ws.Cells.Item(122, 12).Value = 55191
$ws.Cells.Item(122, 13).Value = -182998
$ws.Cells.Item(122, 14).Value = -60091
$ws.Cells.Item(132, 8).Value = 3975.8125
$ws.Cells.Item(132, 9).Value = 4003.5925
$ws.Cells.Item(132, 11).Value = 12010.7775
$ws.Cells.Item(132, 13).Value = -9480.7775
$ws.Cells.Item(135, 8).Value = 5735.643
$ws.Cells.Item(135, 9).Value = 5982.522
$ws.Cells.Item(135, 11).Value = 53842.698
$ws.Cells.Item(135, 13).Value = -51307.698
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 8886.474
$ws.Cells.Item(2, 9).Value = 11941.083
$ws.Cells.Item(2, 10).Value = 3650
$ws.Cells.Item(2, 11).Value = 11941.083
$ws.Cells.Item(2, 12).Value = 3650
$ws.Cells.Item(2, 13).Value = -11828.083
$ws.Cells.Item(2, 14).Value = -3876
$ws.Cells.Item(32, 8).Value = 3038.16
$ws.Cells.Item(32, 9).Value = 3126.261
$ws.Cells.Item(32, 11).Value = 3126.261
$ws.Cells.Item(32, 13).Value = -2839.261
$ws.Cells.Item(116, 8).Value = 8886.474
$ws.Cells.Item(116, 9).Value = 11941.083
$ws.Cells.Item(116, 10).Value = 3650
$ws.Cells.Item(116, 11).Value = 11941.083
$ws.Cells.Item(116, 12).Value = 3650
$ws.Cells.Item(116, 13).Value = -9647.083000000001
$ws.Cells.Item(116, 14).Value = -8238
$ws.Cells.Item(122, 8).Value = 881338.9399999999
$ws.Cells.Item(122, 9).Value = 4643.2
$ws.Cells.Item(122, 11).Value = 13929.6
$ws.Cells.Item(122, 13).Value = -11479.6
$ws.Cells.Item(132, 8).Value = 2948.963
$ws.Cells.Item(132, 9).Value = 2041.5238
$ws.Cells.Item(132, 11).Value = 6124.5714
$ws.Cells.Item(132, 13).Value = -3594.5714
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 8886.474
$ws.Cells.Item(3, 9).Value = 11941.083
$ws.Cells.Item(3, 10).Value = 3650
$ws.Cells.Item(3, 11).Value = 11941.083
$ws.Cells.Item(3, 12).Value = 3650
$ws.Cells.Item(3, 13).Value = -11827.083
$ws.Cells.Item(3, 14).Value = -3878
$ws.Cells.Item(20, 8).Value = 3297.5264
$ws.Cells.Item(20, 9).Value = 2168.818
$ws.Cells.Item(20, 11).Value = 2168.818
$ws.Cells.Item(20, 13).Value = -1921.818
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 791
$ws.Cells.Item(122, 9).Value = 791.6667
$ws.Cells.Item(122, 11).Value = 2375.0001
$ws.Cells.Item(122, 13).Value = 74.9998999999998
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 834576.4399999999
$ws.Cells.Item(4, 9).Value = 870427.6
$ws.Cells.Item(4, 10).Value = 9999
$ws.Cells.Item(4, 11).Value = 2611282.8
$ws.Cells.Item(4, 12).Value = 29997
$ws.Cells.Item(4, 13).Value = -2611170.8
$ws.Cells.Item(4, 14).Value = -30221
$ws.Cells.Item(12, 8).Value = 121.1
$ws.Cells.Item(12, 10).Value = 130.14285
$ws.Cells.Item(12, 12).Value = 390.42855
$ws.Cells.Item(12, 14).Value = -736.4285500000001
$ws.Cells.Item(60, 8).Value = 3693.3333
$ws.Cells.Item(60, 10).Value = 5500
$ws.Cells.Item(60, 12).Value = 16500
$ws.Cells.Item(60, 14).Value = -17002
$ws.Cells.Item(75, 8).Value = 1554
$ws.Cells.Item(75, 10).Value = 2316.3333
$ws.Cells.Item(75, 12).Value = 6948.999899999999
$ws.Cells.Item(75, 14).Value = -8944.999899999999
$ws.Cells.Item(78, 8).Value = 1554
$ws.Cells.Item(78, 10).Value = 2316.3333
$ws.Cells.Item(78, 12).Value = 20846.9997
$ws.Cells.Item(78, 14).Value = -30830.9997
$ws.Cells.Item(129, 8).Value = 2025
$ws.Cells.Item(129, 9).Value = 950
$ws.Cells.Item(129, 10).Value = 3100
$ws.Cells.Item(129, 11).Value = 2850
$ws.Cells.Item(129, 12).Value = 9300
$ws.Cells.Item(129, 13).Value = 2150
$ws.Cells.Item(129, 14).Value = -19300
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 7999.3335
$ws.Cells.Item(70, 9).Value = 0
$ws.Cells.Item(70, 10).Value = 7999.3335
$ws.Cells.Item(70, 11).Value = 0
$ws.Cells.Item(70, 12).Value = 7999.3335
$ws.Cells.Item(70, 13).ClearContents()
$ws.Cells.Item(70, 14).Value = -8539.333500000001
$ws.Cells.Item(73, 8).Value = 7999.3335
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 7999.3335
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 7999.3335
$ws.Cells.Item(73, 13).ClearContents()
$ws.Cells.Item(73, 14).Value = -9871.333500000001
$ws.Cells.Item(80, 8).Value = 2793.1538
$ws.Cells.Item(80, 9).Value = 2998
$ws.Cells.Item(80, 11).Value = 2998
$ws.Cells.Item(80, 13).Value = -2000
$ws.Cells.Item(83, 8).Value = 2793.1538
$ws.Cells.Item(83, 9).Value = 2998
$ws.Cells.Item(83, 11).Value = 14990
$ws.Cells.Item(83, 13).Value = -9998
$ws.Cells.Item(97, 8).Value = 9376.84
$ws.Cells.Item(97, 9).Value = 12168.389
$ws.Cells.Item(97, 11).Value = 12168.389
$ws.Cells.Item(97, 13).Value = -11672.389
$ws.Cells.Item(123, 8).Value = 17062.438
$ws.Cells.Item(123, 10).Value = 17062.438
$ws.Cells.Item(123, 12).Value = 17062.438
$ws.Cells.Item(123, 14).Value = -21962.438
$ws.Cells.Item(132, 8).Value = 3356.6296
$ws.Cells.Item(132, 9).Value = 2363.9048
$ws.Cells.Item(132, 11).Value = 7091.714399999999
$ws.Cells.Item(132, 13).Value = -4561.714399999999
$ws.Cells.Item(141, 8).Value = 100140
$ws.Cells.Item(141, 10).Value = 100140
$ws.Cells.Item(141, 12).Value = 100140
$ws.Cells.Item(141, 14).Value = -110500
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1869.6538
$ws.Cells.Item(22, 9).Value = 1865.3529
$ws.Cells.Item(22, 11).Value = 1865.3529
$ws.Cells.Item(22, 13).Value = -1570.3529
$ws.Cells.Item(27, 8).Value = 1869.6538
$ws.Cells.Item(27, 9).Value = 1865.3529
$ws.Cells.Item(27, 11).Value = 1865.3529
$ws.Cells.Item(27, 13).Value = -1758.3529
$ws.Cells.Item(55, 8).Value = 775.2105
$ws.Cells.Item(55, 9).Value = 737.5599999999999
$ws.Cells.Item(55, 10).Value = 847.61536
$ws.Cells.Item(55, 11).Value = 737.5599999999999
$ws.Cells.Item(55, 12).Value = 847.61536
$ws.Cells.Item(55, 13).Value = -564.5599999999999
$ws.Cells.Item(55, 14).Value = -1193.61536
$ws.Cells.Item(61, 8).Value = 3457.0833
$ws.Cells.Item(61, 9).Value = 3000.2856
$ws.Cells.Item(61, 10).Value = 4096.6
$ws.Cells.Item(61, 11).Value = 3000.2856
$ws.Cells.Item(61, 12).Value = 4096.6
$ws.Cells.Item(61, 13).Value = -2798.2856
$ws.Cells.Item(61, 14).Value = -4500.6
$ws.Cells.Item(113, 8).Value = 3457.0833
$ws.Cells.Item(113, 9).Value = 3000.2856
$ws.Cells.Item(113, 10).Value = 4096.6
$ws.Cells.Item(113, 11).Value = 3000.2856
$ws.Cells.Item(113, 12).Value = 4096.6
$ws.Cells.Item(113, 13).Value = -830.2856000000002
$ws.Cells.Item(113, 14).Value = -8436.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 36500
$ws.Cells.Item(75, 10).Value = 36500
$ws.Cells.Item(75, 12).Value = 36500
$ws.Cells.Item(75, 14).Value = -38372
$ws.Cells.Item(78, 8).Value = 36500
$ws.Cells.Item(78, 10).Value = 36500
$ws.Cells.Item(78, 12).Value = 109500
$ws.Cells.Item(78, 14).Value = -118860
$ws.Cells.Item(100, 8).Value = 24252
$ws.Cells.Item(100, 9).Value = 13002.526
$ws.Cells.Item(100, 10).Value = 77687
$ws.Cells.Item(100, 11).Value = 26005.052
$ws.Cells.Item(100, 12).Value = 155374
$ws.Cells.Item(100, 13).Value = -25464.052
$ws.Cells.Item(100, 14).Value = -156456
$ws.Cells.Item(132, 8).Value = 11880.405
$ws.Cells.Item(132, 9).Value = 14389.173
$ws.Cells.Item(132, 11).Value = 43167.519
$ws.Cells.Item(132, 13).Value = -40637.519
$ws.Cells.Item(136, 8).Value = 6372.25
$ws.Cells.Item(136, 9).Value = 5494
$ws.Cells.Item(136, 11).Value = 16482
$ws.Cells.Item(136, 13).Value = -13932
